$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A3").Value = "Austin"
$ws.Range("B3").Value = "Alison"
$ws.Range("C3").Value = 72

$ws.Range("A4").Select()
